$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.528.85'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").Value = '2.264.66'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '118.49'
$ws.Range("E5").Value = '  +4.36%  '
$ws.Range("D6").Value = '264.91'
$ws.Range("E6").Value = '  -1.02%  '
$ws.Range("D7").Value = '0.643'
$ws.Range("E7").Value = '  +2.67%  '
$ws.Range("D9").Value = '0.617'
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("D10").Value = '47.24'
$ws.Range("E10").Value = '  -2.97%  '
$ws.Range("D11").Value = '0.0940'
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = '9.10'
$ws.Range("E12").Value = '  +2.13%  '
$ws.Range("E13").Value = '  -1.77%  '
$ws.Range("D14").Value = '15.36'
$ws.Range("E14").Value = '  -2.46%  '
$ws.Range("D15").Value = '0.901'
$ws.Range("E15").Value = '  +2.84%  '
$ws.Range("D16").Value = '2.607.50'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").Value = '2.264.01'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").Value = '43.489.97'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").Value = '0.0000109'
$ws.Range("E19").Value = '  +0.51%  '
$ws.Range("D20").Value = '6.83'
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("D21").Value = '71.95'
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").Value = '2.38'
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("D23").Value = '235.06'
$ws.Range("E23").Value = '  +1.02%  '
$ws.Range("D24").Value = '9.46'
$ws.Range("E24").Value = '  -5.01%  '
$ws.Range("E25").Value = '  -0.60%  '
$ws.Range("D27").Value = '11.93'
$ws.Range("E27").Value = '  +3.02%  '
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("E29").Value = '  -0.69%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '171.34'
$ws.Range("E31").Value = '  -1.19%  '
$ws.Range("D32").Value = '21.57'
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '0.129'
$ws.Range("E35").Value = '  +1.56%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '4.18'
$ws.Range("E36").Value = '  +11.13%  '
$ws.Range("D37").Value = '0.0379'
$ws.Range("E37").Value = '  +7.63%  '
$ws.Range("E38").Value = '  -2.33%  '
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("D40").Value = '2.53'
$ws.Range("E40").Value = '  +3.86%  '
$ws.Range("D41").Value = '13.93'
$ws.Range("E41").Value = '  -6.61%  '
$ws.Range("D42").Value = '73.84'
$ws.Range("E42").Value = '  -2.07%  '
$ws.Range("E43").Value = '  -1.09%  '
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("B45").Value = 'THORChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D45").Value = '5.80'
$ws.Range("E45").Value = '  -7.64%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '1.36'
$ws.Range("E46").Value = '  -1.62%  '
$ws.Range("D47").Value = '73.71'
$ws.Range("E47").Value = '  +40.35%  '
$ws.Range("D48").Value = '8.54'
$ws.Range("E48").Value = '  -1.62%  '
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("D50").Value = '0.0996'
$ws.Range("E50").Value = '  -0.28%  '
$ws.Range("D51").Value = '101.07'
$ws.Range("E51").Value = '  -1.14%  '
